# relatorio_financeiro.xlsx - update financial data
# - restyle header row (bold white text on blue fill, thin border, centered)
# - apply currency number format "R$"#,##0.00 to the Valor (R$) column
# - set explicit column widths
# - replace the placeholder "teste2"/"teste3" rows with the real
#   "Total de gastos" / "Saldo final" rows (formulas now point at the
#   single remaining expense row, B4)
# - update a couple of values
# - shrink the used range from A1:C8 to A1:C6 and keep the chart series in sync

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header styling (A1:C1): white bold font, blue fill, thin border, centered ---
$headerCell = $ws.Range("A1")
$headerCell.Font.Color = 16777215        # RGB(255,255,255)
$headerCell.Interior.Color = 12419407    # RGB(79,129,189) -> 004F81BD
$headerCell.Borders.LineStyle = 1
$headerCell.HorizontalAlignment = -4108  # xlCenter
$headerCell.Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Column widths (stored widths: 17 / 13 / 15 characters) ---
$ws.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws.Columns.Item(2).ColumnWidth = 12.166666666666666
$ws.Columns.Item(3).ColumnWidth = 14.166666666666666

# --- Data updates ---
$ws.Range("B3").Value = 300
$ws.Range("B4").Value = 304

$ws.Range("A5").Value = "Total de gastos"
$ws.Range("B5").Formula = "=SUM(B4:B4)"
$ws.Range("A6").Value = "Saldo final"
$ws.Range("B6").Formula = "=B2+B3-B5"

# --- Currency number format for the Valor (R$) column (B2:B6) ---
$ws.Range("B2").NumberFormat = '"R$"#,##0.00'
$ws.Range("B2").Copy()
$ws.Range("B3:B6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Drop the now-obsolete rows 7 and 8 (old totals/balance rows) ---
$ws.Range("A7:C8").Delete()

# --- Keep the chart's series ranges in sync with the shrunk data range ---
$chartObj = $ws.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.Formula = "=SERIES(,'Relatório Financeiro'!`$A`$2:`$A`$6,'Relatório Financeiro'!`$B`$2:`$B`$6,1)"
